$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price (column D) cells to remain plain text,
# matching the original workbook where every D/E cell is stored as text
# (even values like "1.00" or "0.543"). Without this, Excel's COM Value
# setter auto-converts a plain numeric string into a real number.
$numericTextCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D20",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D41",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.615.17"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.735.70"
$ws.Range("E3").Value = "  +18.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "615.43"
$ws.Range("E5").Value = "  +6.46%  "
$ws.Range("D6").Value = "177.15"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "3.732.82"
$ws.Range("E7").Value = "  +18.77%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +10.05%  "
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  +6.94%  "
$ws.Range("D13").Value = "41.14"
$ws.Range("E13").Value = "  +11.17%  "
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "4.360.09"
$ws.Range("E15").Value = "  +18.86%  "
$ws.Range("D16").Value = "3.738.72"
$ws.Range("E16").Value = "  +18.94%  "
$ws.Range("D17").Value = "69.712.28"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +6.69%  "
$ws.Range("D20").Value = "518.61"
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "9.34"
$ws.Range("E22").Value = "  +19.89%  "
$ws.Range("D23").Value = "0.735"
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").Value = "88.78"
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("E25").Value = "  +5.78%  "
$ws.Range("D26").Value = "13.63"
$ws.Range("E26").Value = "  +4.95%  "
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "0.0000126"
$ws.Range("E29").Value = "  +32.35%  "
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +6.41%  "
$ws.Range("E31").Value = "  +8.17%  "
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").Value = "31.57"
$ws.Range("E33").Value = "  +11.63%  "
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +9.26%  "
$ws.Range("E37").Value = "  +8.33%  "
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("D41").Value = "51.37"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("D42").Value = "8.87"
$ws.Range("D43").Value = "44.46"
$ws.Range("E43").Value = "  -8.21%  "
$ws.Range("D44").Value = "425.18"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").Value = "3.063.23"
$ws.Range("E45").Value = "  +9.13%  "
$ws.Range("D46").Value = "2.72"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("D48").Value = "28.01"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  +6.19%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "135.50"
$ws.Range("E51").Value = "  +0.33%  "

# Restore the original (default) style on the cells we reformatted as text,
# so only the NumberFormat-driven text coercion remains and no stray
# direct formatting is left behind on the cell.
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
